$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: before D (so old C stays, new col becomes "px" variant)
# Original layout: C=Primary(mm) D=Lateral(mm) E=Total(mm) F=Leaf(px) G=Lateral_root_count
# Target layout:   C=Primary(px) D=Primary(mm) E=Lateral(px) F=Lateral(mm) G=Total(px) H=Total(mm) I=Leaf(px) J=Lateral_root_count
# So we need to insert a new column before each of C, D, E (shifting things right),
# effectively inserting 3 new columns total at positions C, E, G (old indices), which
# is the same as inserting at D, F, H working left-to-right after each prior insert shifts.
# Simplest: insert a column before C, before E (now F), before G (now I) -- easier to just
# insert three columns one at a time from left to right recomputing positions.

# Insert new column before current C (so old C -> D)
$ws.Columns("C").Insert()
# Insert new column before current E (old D -> E, which is now shifted, new E is old D's spot)
$ws.Columns("E").Insert()
# Insert new column before current G (old E -> ... )
$ws.Columns("G").Insert()

# Now columns are:
# A, B unchanged
# C = new (Primary px)
# D = old C = Primary(mm)
# E = new (Lateral px)
# F = old D = Lateral(mm)
# G = new (Total px)
# H = old E = Total(mm)
# I = old F = Leaf(px)
# J = old G = Lateral_root_count

# Set headers
$ws.Range("C1").Value = "Primary_length(px)"
$ws.Range("D1").Value = "Primary_length(mm)"
$ws.Range("E1").Value = "Lateral_length(px)"
$ws.Range("F1").Value = "Lateral_length(mm)"
$ws.Range("G1").Value = "Total_length(px)"
$ws.Range("H1").Value = "Total_length(mm)"
$ws.Range("I1").Value = "Leaf_size(px)"
$ws.Range("J1").Value = "Lateral_root_count"

# Copy style from B1 (header style) to new header cells C1, E1, G1
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

# Set data values
$data = @{
  2 = @{ C=713.7350647362948; D=29.3959982103616; E=204.509667991878; F=8.422965511036047; G=918.2447327281728; H=37.81896372139765 }
  3 = @{ C=700.9777054234141; D=28.8705717180133; E=145.2670273047587; F=5.982989327073321; G=846.2447327281728; H=34.85356104508662 }
  4 = @{ C=745.4163056034267; D=30.70082649447615; E=139.3675323681472; F=5.740011853824552; G=884.7838379715739; H=36.4408383483007 }
  5 = @{ C=656.1198410471449; D=27.02304906419405; E=12.24264068711928; F=0.5042272147041367; G=668.3624817342642; H=27.52727627889819 }
  6 = @{ C=223.2375900532359; D=9.1942964860707; E=84.52691193458119; F=3.481337928765471; G=307.7645019878171; H=12.67563441483617 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  $ws.Range("C$row").Value = $rowData.C
  $ws.Range("D$row").Value = $rowData.D
  $ws.Range("E$row").Value = $rowData.E
  $ws.Range("F$row").Value = $rowData.F
  $ws.Range("G$row").Value = $rowData.G
  $ws.Range("H$row").Value = $rowData.H
}
